$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$arr = New-Object 'double[,]' 30,3
$arr[0,0] = -2.243617534637451
$arr[0,1] = 1.99995231628418
$arr[0,2] = 0.6877330541610718
$arr[1,0] = -1.472963809967041
$arr[1,1] = 1.918478488922119
$arr[1,2] = 0.357224702835083
$arr[2,0] = -2.048044681549072
$arr[2,1] = 1.973574161529541
$arr[2,2] = 0.7422993183135986
$arr[3,0] = -1.858012676239014
$arr[3,1] = 2.016510963439941
$arr[3,2] = 0.136701762676239
$arr[4,0] = -1.831352233886719
$arr[4,1] = 1.850147724151612
$arr[4,2] = 0.5375880002975464
$arr[5,0] = -1.778313636779785
$arr[5,1] = 1.832831382751465
$arr[5,2] = 0.6446512937545776
$arr[6,0] = -1.913262844085693
$arr[6,1] = 1.87591552734375
$arr[6,2] = 0.3856300711631775
$arr[7,0] = -2.06378698348999
$arr[7,1] = 1.738418102264404
$arr[7,2] = 0.8020429015159607
$arr[8,0] = -1.788212299346924
$arr[8,1] = 1.889766216278076
$arr[8,2] = 0.3115454912185669
$arr[9,0] = -1.937233924865723
$arr[9,1] = 1.803182125091553
$arr[9,2] = 0.8862782716751099
$arr[10,0] = -2.046533107757568
$arr[10,1] = 1.845665454864502
$arr[10,2] = 0.7709857225418091
$arr[11,0] = -1.904958724975586
$arr[11,1] = 1.821238040924072
$arr[11,2] = 0.7353460788726807
$arr[12,0] = -1.893377780914306
$arr[12,1] = 1.849058628082276
$arr[12,2] = 0.4389755129814148
$arr[13,0] = -1.822815418243408
$arr[13,1] = 1.71094799041748
$arr[13,2] = 0.5070880651473999
$arr[14,0] = -1.894952774047852
$arr[14,1] = 1.727813243865967
$arr[14,2] = 0.551978588104248
$arr[15,0] = -1.988187789916992
$arr[15,1] = 1.761092185974121
$arr[15,2] = 0.5642168521881104
$arr[16,0] = -1.921907901763916
$arr[16,1] = 1.732572555541992
$arr[16,2] = 0.5020300149917603
$arr[17,0] = -1.923232555389404
$arr[17,1] = 1.782273769378662
$arr[17,2] = 0.3185268342494964
$arr[18,0] = -1.910520076751709
$arr[18,1] = 1.814604759216309
$arr[18,2] = 0.3955313861370086
$arr[19,0] = -2.009746074676514
$arr[19,1] = 1.783975601196289
$arr[19,2] = 0.4257155656814575
$arr[20,0] = -1.97176456451416
$arr[20,1] = 1.745009422302246
$arr[20,2] = 0.4838592410087585
$arr[21,0] = -1.845728397369385
$arr[21,1] = 1.672563552856445
$arr[21,2] = 0.5211508870124817
$arr[22,0] = -2.323234558105469
$arr[22,1] = 1.572612762451172
$arr[22,2] = 0.7384862899780273
$arr[23,0] = -3.269093036651612
$arr[23,1] = 0.9374399185180664
$arr[23,2] = 1.780844926834106
$arr[24,0] = -1.246297359466553
$arr[24,1] = 0.5547242164611816
$arr[24,2] = 1.987784385681152
$arr[25,0] = -2.729167938232422
$arr[25,1] = 0.9446659088134766
$arr[25,2] = 2.494723320007324
$arr[26,0] = -2.581250190734864
$arr[26,1] = -1.341118812561035
$arr[26,2] = 3.135838031768799
$arr[27,0] = -3.531517744064331
$arr[27,1] = -1.318747520446777
$arr[27,2] = 3.353370428085327
$arr[28,0] = -4.357868194580078
$arr[28,1] = -0.8888802528381348
$arr[28,2] = 1.668703079223633
$arr[29,0] = -2.958048820495605
$arr[29,1] = 1.071966171264648
$arr[29,2] = 2.359493732452393

$ws.Range("A2:C31").Value = $arr

Write-Host "Done"